$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert a new column before E ("Surprise"), pushing Phone/Column with
# int/Empty column one column to the right (E->F, F->G, G->H). ---
$ws.Columns("E:E").Insert()
$ws.Range("E1").Value = "Surprise"

# --- Insert a new, entirely blank row before the old row 3 (John Smith),
# pushing that record down to row 4. ---
$ws.Rows("3:3").Insert()

# --- New values typed in a far-off column I for the two data rows. ---
$ws.Range("I2").Value = "dsdsqd"
$ws.Range("I4").Value = "qdqsdq"

# --- The row insert does not relocate the existing hyperlink anchored on
# the old row 3 (jean.lefebvre@mail.com), so rebuild it pointing at C4. ---
$ws.Range("C2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:john.smith@mail.com", [Type]::Missing, [Type]::Missing, "john.smith@mail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:jean.lefebvre@mail.com", [Type]::Missing, [Type]::Missing, "jean.lefebvre@mail.com")

# Adding a hyperlink re-styles the cell with the built-in "Hyperlink" look;
# put the original font back so C2/C4 keep their pre-existing appearance.
$ws.Range("C4").Font.Name = "Droid Sans"
$ws.Range("C4").Font.Color = 16711680
$ws.Range("C4").Font.Underline = -4142
$ws.Range("C4").Font.Size = 10
$ws.Range("C2").Font.Name = "Droid Sans"
$ws.Range("C2").Font.Color = 16711680
$ws.Range("C2").Font.Underline = -4142
$ws.Range("C2").Font.Size = 10
$wb.Styles.Item("Hyperlink").Delete()

# --- Restore the selection left behind on the sheet. ---
$ws.Range("C11").Select()

Write-Host "edit applied"
